# bank_acc logic in print_docs func is added
# Applies the six textual field updates described by the diff.

$d = $word.ActiveDocument

# 1) Vessel name: "АКАДЕМИК МСТИСЛАВ КЕЛДЫШ" -> "28 MAY"
$d.Content.Find.Execute(
    "АКАДЕМИК МСТИСЛАВ КЕЛДЫШ", $true, $false, $false, $false, $false,
    $true, 1, $false, "28 MAY", 2)

# 2) Registration number: 780270 -> 120378
$d.Content.Find.Execute(
    "780270", $true, $false, $false, $false, $false,
    $true, 1, $false, "120378", 2)

# 3) IMO number: 7811018 -> "--"
$d.Content.Find.Execute(
    "7811018", $true, $false, $false, $false, $false,
    $true, 1, $false, "--", 2)

# 4) Authority basis: "Доверенности № 123456 от 02.02.2024" -> "Устава"
$d.Content.Find.Execute(
    "Доверенности № 123456 от 02.02.2024", $true, $false, $false, $false, $false,
    $true, 1, $false, "Устава", 2)

# 5) Settlement currency: "РУБ" (underlined) -> "EUR"
$d.Content.Find.Execute(
    "РУБ", $true, $false, $false, $false, $false,
    $true, 1, $false, "EUR", 2)

# 6) Bank account details block
$d.Content.Find.Execute(
    "р/с 301018101000000000634, калининградский филиал ПАО АКБ «Связь-Банк», Калининград, БИК 042748898 к/с 30101810700000000898",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "IBAN 3532465667668, DANSKE BANK A/S Lietuvos filialas, BIC SMPOLT22XXX", 2)
